{"js": "// Add a new \"Segunda Quest\u00e3o:\" paragraph at the end of the document body,\n// mirroring the structure/formatting of the existing \"Primeira Quest\u00e3o:\"\n// paragraph: a bold label run followed by normal (non-bold) runs holding\n// the answer text.\n\nconst body = context.document.body;\n\n// Insert a brand-new, empty paragraph at the very end of the document.\nconst newPara = body.insertParagraph(\"\", Word.InsertLocation.end);\nnewPara.alignment = Word.Alignment.justified;\n\n// Bold label run: \"Segunda Quest\u00e3o:\"\nconst label = newPara.insertText(\"Segunda Quest\u00e3o:\", Word.InsertLocation.end);\nlabel.font.bold = true;\n\n// Remaining answer text, split into several runs (mirrors how Word\n// naturally splits runs while typing/auto-correcting); each one is\n// explicitly reset to non-bold so it doesn't inherit the label's bold\n// formatting.\nconst parts = [\n  \" \",\n  \"Sim o m\u00e9todo \u00e1gil combina com esse modelo,\",\n  \" \",\n  \"pois,\",\n  \" esse \",\n  \"tipo de  modelo combina com essa abordagem\"\n];\n\nfor (const part of parts) {\n  const run = newPara.insertText(part, Word.InsertLocation.end);\n  run.font.bold = false;\n}\n\nawait context.sync();\n", "ps1": "# Add a new \"Segunda Quest\u00e3o:\" paragraph at the end of the document,\n# mirroring the structure/formatting of the existing \"Primeira Quest\u00e3o:\"\n# paragraph: bold label run followed by normal (non-bold) runs for the\n# answer text.\n\n$doc = $word.ActiveDocument\n\n# Move to the end of the last paragraph and insert a new paragraph mark.\n$lastPara = $doc.Paragraphs.Item($doc.Paragraphs.Count)\n$tail = $lastPara.Range\n$tail.Collapse(0)            # wdCollapseEnd\n$tail.InsertParagraphAfter()\n\n# The newly created paragraph is now the last one in the document.\n$newPara = $doc.Paragraphs.Item($doc.Paragraphs.Count)\n$newPara.Format.Alignment = 3   # wdAlignParagraphJustify\n\n$r = $newPara.Range\n$r.Collapse(0)\n\n# Bold label run: \"Segunda Quest\u00e3o:\"\n$r.InsertAfter(\"Segunda Quest\u00e3o:\")\n$r.Font.Bold = 1\n\n# Remaining answer text, split into several runs (mirrors how Word\n# naturally splits runs while typing/auto-correcting); each one explicitly\n# reset to non-bold so it doesn't inherit the label's bold formatting.\n$parts = @(\n  \" \",\n  \"Sim o m\u00e9todo \u00e1gil combina com esse modelo,\",\n  \" \",\n  \"pois,\",\n  \" esse \",\n  \"tipo de  modelo combina com essa abordagem\"\n)\n\nforeach ($part in $parts) {\n  $r.Collapse(0)\n  $r.InsertAfter($part)\n  $r.Font.Bold = 0\n}\n"}
